$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2005"
$ws.Range("B2").Value = "**神戸大学** <br> [市販ミンチ肉における黄色ブドウ球菌汚染調査と分離株の性状](https://www.jstage.jst.go.jp/article/jsfm1994/23/4/23_4_217/_pdf/-char/ja) <br> （日本食品微生物学会雑誌, 23 (4), 217-222, 2006）"
$ws.Range("C2").Value = "済"

# Update row 3
$ws.Range("A3").Value = "2000-2002"
$ws.Range("B3").Value = "**東京都健康安全研究センター** <br> [市販生鮮青果物の食品細菌学的調査](https://www.tmiph.metro.tokyo.lg.jp/files/archive/issue/kenkyunenpo/nenpou55/55-21.pdf) <br>（東京健康安全研究センター年報, 55, 2004）"
$ws.Range("C3").Value = "済"

# Remove row 4 entirely (was: 2006, 和歌山県立医科大学 reference, 未登録)
$ws.Rows.Item(4).Delete()
